$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.719.26'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.853.81'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.06'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6400'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07484'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2985'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.38'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07646'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.854.09'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.043'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6879'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.81'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009493'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +5.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.053'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.720.00'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.101.86'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '235.65'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.63'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.410'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.02%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.37'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.485'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06278'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.493'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.277'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +6.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.147'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.089'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.902'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.170'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7287'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.605'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.843'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01784'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.202.25'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9230'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.145'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.011.80'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.93'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '66.01'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.82%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000119'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.207'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4061'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05800'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.93%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.648'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.67%  '
